# "Renouveau DonneeTest et Changement commentaire"
# Regenerate the test data in the "Resultat" sheet: new random-looking values
# for row 2 (A2:I2, skipping D2) and column C (C3:C80), then drop the now
# unused rows 81:93 so the sheet shrinks from A1:I93 to A1:I80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (multi-column summary row) ---
$ws.Range("A2").Value = 79
$ws.Range("B2").Value = 121
$ws.Range("C2").Value = 40
$ws.Range("E2").Value = 40
$ws.Range("F2").Value = 41
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 45
$ws.Range("I2").Value = 5

# --- Column C data series (rows 3 through 80) ---
$colC = @(112,83,78,102,76,28,48,80,77,41,89,15,23,26,39,67,52,90,62,118,55,95,75,80,3,114,0,40,69,69,95,55,97,4,34,69,107,8,79,13,42,105,58,116,13,17,69,36,26,8,53,39,82,16,11,121,41,97,28,67,118,71,50,62,13,64,116,72,112,42,105,21,5,116,39,116,90,44)

for ($i = 0; $i -lt $colC.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}

# --- Remove the now-obsolete tail of the data series (rows 81:93) ---
$ws.Range("A81:A93").EntireRow.Delete()
